# Add yearly repeated measures: change value type of "height_" from "integer" to "decimal"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")

$ws.Range("B4").Value = "decimal"

# Update the active selection to B4 as shown in the saved file
$ws.Range("B4").Select()
